$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 4, "301.42"),
    @(2, 5, "0.32%"),
    @(2, 7, "8"),
    @(3, 4, "32.79"),
    @(3, 5, "4.29%"),
    @(3, 7, "8"),
    @(4, 4, "4.940"),
    @(4, 5, "-2.82%"),
    @(4, 7, "8"),
    @(5, 4, "0.07764"),
    @(5, 5, "-1.30%"),
    @(5, 7, "8"),
    @(6, 4, "1.946"),
    @(6, 5, "-14.04%"),
    @(6, 7, "8"),
    @(7, 4, "7.849"),
    @(7, 5, "0.51%"),
    @(7, 7, "8"),
    @(8, 4, "3.803"),
    @(8, 5, "-0.87%"),
    @(8, 7, "8"),
    @(9, 4, "0.9197"),
    @(9, 5, "-0.13%"),
    @(9, 7, "8"),
    @(10, 4, "0.1775"),
    @(10, 5, "1.78%"),
    @(10, 7, "8"),
    @(11, 4, "0.07907"),
    @(11, 5, "4.14%"),
    @(11, 7, "8"),
    @(12, 4, "0.08615"),
    @(12, 5, "-5.86%"),
    @(12, 7, "8"),
    @(13, 4, "0.03154"),
    @(13, 5, "5.03%"),
    @(13, 7, "8"),
    @(14, 4, "0.1004"),
    @(14, 5, "0.12%"),
    @(14, 7, "8"),
    @(15, 4, "0.001514"),
    @(15, 5, "0.43%"),
    @(15, 7, "8"),
    @(16, 4, "0.005860"),
    @(16, 5, "-2.54%"),
    @(16, 7, "8"),
    @(17, 7, "8"),
    @(18, 4, "2.154"),
    @(18, 5, "-3.99%"),
    @(18, 7, "8"),
    @(19, 4, "0.3341"),
    @(19, 5, "2.13%"),
    @(19, 7, "8"),
    @(20, 5, "2.13%"),
    @(20, 7, "8"),
    @(21, 4, "4.309"),
    @(21, 5, "7.99%"),
    @(21, 7, "8"),
    @(22, 4, "0.1992"),
    @(22, 5, "16.54%"),
    @(22, 7, "8"),
    @(23, 4, "0.04563"),
    @(23, 5, "-1.12%"),
    @(23, 7, "8"),
    @(24, 4, "0.001227"),
    @(24, 5, "-2.11%"),
    @(24, 7, "8"),
    @(25, 4, "0.004422"),
    @(25, 5, "-1.06%"),
    @(25, 7, "8"),
    @(26, 4, "0.0001251"),
    @(26, 5, "0.25%"),
    @(26, 7, "8"),
    @(27, 7, "8"),
    @(28, 7, "8"),
    @(29, 7, "8"),
    @(30, 7, "8"),
    @(31, 7, "8"),
    @(32, 7, "8"),
    @(33, 7, "8"),
    @(34, 7, "8"),
    @(35, 7, "8"),
    @(36, 7, "8"),
    @(37, 7, "8"),
    @(38, 7, "8"),
    @(39, 4, "0.01702"),
    @(39, 5, "-2.27%"),
    @(39, 7, "8"),
    @(40, 4, "0.04718"),
    @(40, 5, "1.92%"),
    @(40, 7, "8"),
    @(41, 4, "0.007464"),
    @(41, 5, "5.59%"),
    @(41, 7, "8"),
    @(42, 5, "-0.27%"),
    @(42, 7, "8"),
    @(43, 4, "0.002342"),
    @(43, 5, "7.12%"),
    @(43, 7, "8"),
    @(44, 4, "0.01044"),
    @(44, 5, "7.07%"),
    @(44, 7, "8"),
    @(45, 4, "0.00006240"),
    @(45, 5, "-0.45%"),
    @(45, 7, "8"),
    @(46, 4, "0.00000000751"),
    @(46, 5, "0.27%"),
    @(46, 7, "8"),
    @(47, 4, "0.8204"),
    @(47, 5, "-28.87%"),
    @(47, 7, "8"),
    @(48, 4, "0.003103"),
    @(48, 5, "-61.15%"),
    @(48, 7, "8"),
    @(49, 4, "0.00002102"),
    @(49, 5, "0.27%"),
    @(49, 7, "8"),
    @(50, 4, "0.0002002"),
    @(50, 5, "0.27%"),
    @(50, 7, "8"),
    @(51, 7, "8")
)

foreach ($item in $data) {
    $r = $item[0]
    $c = $item[1]
    $val = $item[2]
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}
Write-Output "Applied $($data.Count) cell updates"